$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 22:18:19"
$ws.Range("N2").Value = "-2.9 °C 21:55 TU"
$ws.Range("O2").Value = "0.3 °C"
$ws.Range("E3").Value = "2026-02-20 22:18:21"
$ws.Range("E4").Value = "2026-02-20 22:18:24"
$ws.Range("H4").Value = "'60%"
$ws.Range("J4").Value = "1023.0 hPa"
$ws.Range("N4").Value = "3.9 °C 21:52 TU"
$ws.Range("O4").Value = "9.8 °C"
$ws.Range("E5").Value = "2026-02-20 22:18:26"
$ws.Range("E6").Value = "2026-02-20 22:18:29"
$ws.Range("J6").Value = "1022.9 hPa"
$ws.Range("E7").Value = "2026-02-20 22:18:31"
$ws.Range("H7").Value = "'49%"
$ws.Range("J7").Value = "1022.8 hPa"
$ws.Range("E8").Value = "2026-02-20 22:18:34"
$ws.Range("J8").Value = "1023.1 hPa"
$ws.Range("E9").Value = "2026-02-20 22:18:36"
$ws.Range("E10").Value = "2026-02-20 22:18:39"
$ws.Range("H10").Value = "'80%"
$ws.Range("O10").Value = "7.4 °C"
$ws.Range("E11").Value = "2026-02-20 22:18:41"
$ws.Range("H11").Value = "'34%"
$ws.Range("E12").Value = "2026-02-20 22:18:44"
$ws.Range("E13").Value = "2026-02-20 22:18:46"
$ws.Range("H13").Value = "'46%"
$ws.Range("J13").Value = "1024.2 hPa"
$ws.Range("N13").Value = "-0.7 °C 21:51 TU"
$ws.Range("O13").Value = "6.0 °C"
$ws.Range("E14").Value = "2026-02-20 22:18:48"
$ws.Range("H14").Value = "'60%"
$ws.Range("E15").Value = "2026-02-20 22:18:51"
$ws.Range("O15").Value = "13.2 °C"
$ws.Range("E16").Value = "2026-02-20 22:18:53"
$ws.Range("O16").Value = "-3.0 °C"
$ws.Range("E17").Value = "2026-02-20 22:18:56"
$ws.Range("E18").Value = "2026-02-20 22:18:58"
$ws.Range("J18").Value = "1023.2 hPa"
$ws.Range("O18").Value = "7.7 °C"
$ws.Range("E19").Value = "2026-02-20 22:19:01"
$ws.Range("E20").Value = "2026-02-20 22:19:03"
$ws.Range("E21").Value = "2026-02-20 22:19:06"
$ws.Range("H21").Value = "'38%"
$ws.Range("J21").Value = "1023.1 hPa"
$ws.Range("O21").Value = "9.0 °C"
$ws.Range("E22").Value = "2026-02-20 22:19:08"
$ws.Range("H22").Value = "'48%"
$ws.Range("E23").Value = "2026-02-20 22:19:11"
$ws.Range("E24").Value = "2026-02-20 22:19:13"
$ws.Range("H24").Value = "'68%"
$ws.Range("J24").Value = "1025.6 hPa"
$ws.Range("O24").Value = "9.4 °C"
$ws.Range("E25").Value = "2026-02-20 22:19:15"
$ws.Range("E26").Value = "2026-02-20 22:19:18"
$ws.Range("J26").Value = "1022.1 hPa"
$ws.Range("E27").Value = "2026-02-20 22:19:21"
$ws.Range("E28").Value = "2026-02-20 22:19:23"
$ws.Range("J28").Value = "1023.4 hPa"
$ws.Range("E29").Value = "2026-02-20 22:19:26"
$ws.Range("E30").Value = "2026-02-20 22:19:28"
$ws.Range("H30").Value = "'60%"
$ws.Range("J30").Value = "1022.7 hPa"
$ws.Range("E31").Value = "2026-02-20 22:19:30"
$ws.Range("J31").Value = "1021.9 hPa"
$ws.Range("K31").Value = "13.2 MJ/m2"
$ws.Range("E32").Value = "2026-02-20 22:19:33"
$ws.Range("E33").Value = "2026-02-20 22:19:35"
$ws.Range("H33").Value = "'43%"
$ws.Range("J33").Value = "1023.5 hPa"
$ws.Range("O33").Value = "5.8 °C"
$ws.Range("E34").Value = "2026-02-20 22:19:38"
$ws.Range("E35").Value = "2026-02-20 22:19:40"
$ws.Range("J35").Value = "1027.0 hPa"
$ws.Range("E36").Value = "2026-02-20 22:19:43"
$ws.Range("J36").Value = "1022.8 hPa"
$ws.Range("E37").Value = "2026-02-20 22:19:45"
$ws.Range("J37").Value = "1024.9 hPa"
$ws.Range("O37").Value = "4.5 °C"
$ws.Range("E38").Value = "2026-02-20 22:19:48"
$ws.Range("H38").Value = "'69%"
$ws.Range("O38").Value = "8.7 °C"
$ws.Range("E39").Value = "2026-02-20 22:19:50"
$ws.Range("E40").Value = "2026-02-20 22:19:52"
$ws.Range("H40").Value = "'38%"
$ws.Range("J40").Value = "1023.9 hPa"
$ws.Range("N40").Value = "3.4 °C 21:52 TU"
$ws.Range("O40").Value = "10.1 °C"
$ws.Range("E41").Value = "2026-02-20 22:19:55"
$ws.Range("H41").Value = "'49%"
$ws.Range("J41").Value = "1023.3 hPa"
$ws.Range("E42").Value = "2026-02-20 22:19:57"
$ws.Range("O42").Value = "10.5 °C"
$ws.Range("E43").Value = "2026-02-20 22:20:00"
$ws.Range("H43").Value = "'77%"
$ws.Range("O43").Value = "4.9 °C"
$ws.Range("E44").Value = "2026-02-20 22:20:02"
$ws.Range("M44").Value = "-0.8 °C 21:58 TU"
$ws.Range("O44").Value = "-4.4 °C"
$ws.Range("E45").Value = "2026-02-20 22:20:04"
$ws.Range("J45").Value = "1029.9 hPa"
$ws.Range("E46").Value = "2026-02-20 22:20:07"
$ws.Range("J46").Value = "1026.6 hPa"
